$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value2 = "64.150.42"
$ws.Range("E2").Value2 = "  +4.41%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value2 = "2.755.67"
$ws.Range("E3").Value2 = "  +3.71%  "

# Row 4
$ws.Range("E4").Value2 = "  +0.13%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "580.32"
$ws.Range("E5").Value2 = "  -0.23%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "157.67"
$ws.Range("E6").Value2 = "  +8.21%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value2 = "0.629"
$ws.Range("E7").Value2 = "  +4.93%  "

# Row 8
$ws.Range("E8").Value2 = "  +0.27%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = "2.755.11"
$ws.Range("E9").Value2 = "  +3.18%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "6.77"
$ws.Range("E10").Value2 = "  +2.58%  "

# Row 11
$ws.Range("E11").Value2 = "  +2.37%  "

# Row 12
$ws.Range("E12").Value2 = "  +3.84%  "

# Row 13
$ws.Range("E13").Value2 = "  +0.68%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = "3.240.40"
$ws.Range("E14").Value2 = "  +3.86%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = "27.28"
$ws.Range("E15").Value2 = "  +3.19%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = "64.071.27"
$ws.Range("E16").Value2 = "  +4.44%  "

# Row 17
$ws.Range("E17").Value2 = "  +6.17%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = "2.754.97"
$ws.Range("E18").Value2 = "  +3.42%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = "12.10"
$ws.Range("E19").Value2 = "  +3.47%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "4.96"
$ws.Range("E20").Value2 = "  +3.76%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "364.11"
$ws.Range("E21").Value2 = "  +2.45%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "6.99"
$ws.Range("E22").Value2 = "  +1.09%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "0.544"
$ws.Range("E23").Value2 = "  +2.87%  "

# Row 24
$ws.Range("E24").Value2 = "  -0.17%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "67.09"
$ws.Range("E25").Value2 = "  +4.63%  "

# Row 26
$ws.Range("E26").Value2 = "  +5.81%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "8.65"
$ws.Range("E27").Value2 = "  +1.91%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = "0.999"
$ws.Range("E28").Value2 = "  +0.11%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = "0.0₃0919"
$ws.Range("E29").Value2 = "  +11.47%  "

# Row 30
$ws.Range("E30").Value2 = "  +0.35%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "7.22"
$ws.Range("E31").Value2 = "  +6.46%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = "1.26"
$ws.Range("E32").Value2 = "  +13.45%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = "173.97"
$ws.Range("E33").Value2 = "  +3.88%  "

# Row 34
$ws.Range("E34").Value2 = "  +0.15%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value2 = "20.63"
$ws.Range("E35").Value2 = "  +2.69%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = "4.93"
$ws.Range("E36").Value2 = "  +3.94%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = "1.46"
$ws.Range("E37").Value2 = "  +8.81%  "

# Row 38
$ws.Range("E38").Value2 = "  +5.68%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = "1.01"
$ws.Range("E39").Value2 = "  +11.14%  "

# Row 40
$ws.Range("E40").Value2 = "  +3.81%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = "340.84"
$ws.Range("E41").Value2 = "  -0.84%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = "6.19"
$ws.Range("E42").Value2 = "  +15.30%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "39.50"
$ws.Range("E43").Value2 = "  +2.53%  "

# Row 44
$ws.Range("B44").Value2 = "InjectiveProtocol"
$ws.Range("C44").Value2 = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = "22.38"
$ws.Range("E44").Value2 = "  +5.59%  "

# Row 45
$ws.Range("B45").Value2 = "EnergySwap"
$ws.Range("C45").Value2 = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = "21.83"
$ws.Range("E45").Value2 = "  +6.63%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "0.0600"
$ws.Range("E46").Value2 = "  +3.29%  "

# Row 47
$ws.Range("B47").Value2 = "VeChain"
$ws.Range("C47").Value2 = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "0.0261"
$ws.Range("E47").Value2 = "  +3.73%  "

# Row 48
$ws.Range("B48").Value2 = "Mantle"
$ws.Range("C48").Value2 = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "0.648"
$ws.Range("E48").Value2 = "  +3.13%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = "137.59"
$ws.Range("E49").Value2 = "  +0.47%  "

# Row 50
$ws.Range("E50").Value2 = "  +2.72%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = "1.00"
$ws.Range("E51").Value2 = "  +0.64%  "
